$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These cells hold numeric-looking figures stored as text in the workbook.
# Prefix with an apostrophe so Excel keeps storing them as text (matching
# the original shared-string representation) instead of coercing to a
# number.

# Enterprises density (per 1000 people) -- row 12
$ws.Range("B12").Value = "'36.48"
$ws.Range("C12").Value = "'11.28"
$ws.Range("D12").Value = "'47.76"

# Employment (% of total) -- row 13
$ws.Range("B13").Value = "'6.83"
$ws.Range("C13").Value = "'39.18"
$ws.Range("D13").Value = "'46.01"

# Enterprises (% of total) -- row 15
$ws.Range("B15").Value = "'75.24"
$ws.Range("C15").Value = "'23.26"
